$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2288423333333333
$ws.Range("H2").Value = 0.686527
$ws.Range("I2").Value = 0.001039481910007428
$ws.Range("J2").Value = 0.001039481910007428
$ws.Range("M2").Value = 0.2288423333333333
$ws.Range("N2").Value = 0.686527
$ws.Range("O2").Value = 0.001039481910007428
$ws.Range("P2").Value = 0.001039481910007428
$ws.Range("Q2").Value = 0.05236881352544445
$ws.Range("R2").Value = 0.471319321729
$ws.Range("S2").Value = 0.00000108052264123269
$ws.Range("T2").Value = 0.00000108052264123269
$ws.Range("G3").Value = 0.2288423333333333
$ws.Range("H3").Value = 0.686527
$ws.Range("I3").Value = 0.001039481910007428
$ws.Range("J3").Value = 0.001039481910007428
$ws.Range("M3").Value = 132.0967866666666
$ws.Range("N3").Value = 396.29036
$ws.Range("O3").Value = 0.6000298026593724
$ws.Range("P3").Value = 0.6000298026593726
$ws.Range("Q3").Value = 30.22933688663555
$ws.Range("R3").Value = 272.06403197972
$ws.Range("S3").Value = 0.0006237201253297443
$ws.Range("T3").Value = 0.0006237201253297444
$ws.Range("G4").Value = 0.2288423333333333
$ws.Range("H4").Value = 0.686527
$ws.Range("I4").Value = 0.001039481910007428
$ws.Range("J4").Value = 0.001039481910007428
$ws.Range("O4").Value = 0.0003198253300583355
$ws.Range("P4").Value = 0.0003198253300583356
$ws.Range("Q4").Value = 0.01611271240922222
$ws.Range("R4").Value = 0.145014411683
$ws.Range("S4").Value = 0.0000003324526449577945
$ws.Range("T4").Value = 0.0000003324526449577946
$ws.Range("G5").Value = 0.2288423333333333
$ws.Range("H5").Value = 0.686527
$ws.Range("I5").Value = 0.001039481910007428
$ws.Range("J5").Value = 0.001039481910007428
$ws.Range("M5").Value = 87.75433733333334
$ws.Range("N5").Value = 263.263012
$ws.Range("O5").Value = 0.3986108901005617
$ws.Range("P5").Value = 0.3986108901005617
$ws.Range("Q5").Value = 20.08190731548045
$ws.Range("R5").Value = 180.737165839324
$ws.Range("S5").Value = 0.0004143488093914927
$ws.Range("T5").Value = 0.0004143488093914927
$ws.Range("G6").Value = 132.0967866666666
$ws.Range("H6").Value = 396.29036
$ws.Range("I6").Value = 0.6000298026593724
$ws.Range("J6").Value = 0.6000298026593726
$ws.Range("M6").Value = 0.2288423333333333
$ws.Range("N6").Value = 0.686527
$ws.Range("O6").Value = 0.001039481910007428
$ws.Range("P6").Value = 0.001039481910007428
$ws.Range("Q6").Value = 30.22933688663555
$ws.Range("R6").Value = 272.06403197972
$ws.Range("S6").Value = 0.0006237201253297443
$ws.Range("T6").Value = 0.0006237201253297444
$ws.Range("G7").Value = 132.0967866666666
$ws.Range("H7").Value = 396.29036
$ws.Range("I7").Value = 0.6000298026593724
$ws.Range("J7").Value = 0.6000298026593726
$ws.Range("M7").Value = 132.0967866666666
$ws.Range("N7").Value = 396.29036
$ws.Range("O7").Value = 0.6000298026593724
$ws.Range("P7").Value = 0.6000298026593726
$ws.Range("Q7").Value = 17449.56104765884
$ws.Range("R7").Value = 157046.0494289296
$ws.Range("S7").Value = 0.3600357640794454
$ws.Range("T7").Value = 0.3600357640794456
$ws.Range("G8").Value = 132.0967866666666
$ws.Range("H8").Value = 396.29036
$ws.Range("I8").Value = 0.6000298026593724
$ws.Range("J8").Value = 0.6000298026593726
$ws.Range("O8").Value = 0.0003198253300583355
$ws.Range("P8").Value = 0.0003198253300583356
$ws.Range("Q8").Value = 9.300890716937776
$ws.Range("R8").Value = 83.70801645243999
$ws.Range("S8").Value = 0.0001919047296803717
$ws.Range("T8").Value = 0.0001919047296803718
$ws.Range("G9").Value = 132.0967866666666
$ws.Range("H9").Value = 396.29036
$ws.Range("I9").Value = 0.6000298026593724
$ws.Range("J9").Value = 0.6000298026593726
$ws.Range("M9").Value = 87.75433733333334
$ws.Range("N9").Value = 263.263012
$ws.Range("O9").Value = 0.3986108901005617
$ws.Range("P9").Value = 0.3986108901005617
$ws.Range("Q9").Value = 11592.06597779603
$ws.Range("R9").Value = 104328.5938001643
$ws.Range("S9").Value = 0.2391784137249169
$ws.Range("T9").Value = 0.2391784137249169
$ws.Range("I10").Value = 0.0003198253300583355
$ws.Range("J10").Value = 0.0003198253300583356
$ws.Range("M10").Value = 0.2288423333333333
$ws.Range("N10").Value = 0.686527
$ws.Range("O10").Value = 0.001039481910007428
$ws.Range("P10").Value = 0.001039481910007428
$ws.Range("Q10").Value = 0.01611271240922222
$ws.Range("R10").Value = 0.145014411683
$ws.Range("S10").Value = 0.0000003324526449577945
$ws.Range("T10").Value = 0.0000003324526449577946
$ws.Range("I11").Value = 0.0003198253300583355
$ws.Range("J11").Value = 0.0003198253300583356
$ws.Range("M11").Value = 132.0967866666666
$ws.Range("N11").Value = 396.29036
$ws.Range("O11").Value = 0.6000298026593724
$ws.Range("P11").Value = 0.6000298026593726
$ws.Range("Q11").Value = 9.300890716937776
$ws.Range("R11").Value = 83.70801645243999
$ws.Range("S11").Value = 0.0001919047296803717
$ws.Range("T11").Value = 0.0001919047296803718
$ws.Range("I12").Value = 0.0003198253300583355
$ws.Range("J12").Value = 0.0003198253300583356
$ws.Range("O12").Value = 0.0003198253300583355
$ws.Range("P12").Value = 0.0003198253300583356
$ws.Range("S12").Value = 0.0000001022882417469232
$ws.Range("T12").Value = 0.0000001022882417469233
$ws.Range("I13").Value = 0.0003198253300583355
$ws.Range("J13").Value = 0.0003198253300583356
$ws.Range("M13").Value = 87.75433733333334
$ws.Range("N13").Value = 263.263012
$ws.Range("O13").Value = 0.3986108901005617
$ws.Range("P13").Value = 0.3986108901005617
$ws.Range("Q13").Value = 6.178753640194222
$ws.Range("R13").Value = 55.608782761748
$ws.Range("S13").Value = 0.0001274858594912591
$ws.Range("T13").Value = 0.0001274858594912591
$ws.Range("G14").Value = 87.75433733333334
$ws.Range("H14").Value = 263.263012
$ws.Range("I14").Value = 0.3986108901005617
$ws.Range("J14").Value = 0.3986108901005617
$ws.Range("M14").Value = 0.2288423333333333
$ws.Range("N14").Value = 0.686527
$ws.Range("O14").Value = 0.001039481910007428
$ws.Range("P14").Value = 0.001039481910007428
$ws.Range("Q14").Value = 20.08190731548045
$ws.Range("R14").Value = 180.737165839324
$ws.Range("S14").Value = 0.0004143488093914927
$ws.Range("T14").Value = 0.0004143488093914927
$ws.Range("G15").Value = 87.75433733333334
$ws.Range("H15").Value = 263.263012
$ws.Range("I15").Value = 0.3986108901005617
$ws.Range("J15").Value = 0.3986108901005617
$ws.Range("M15").Value = 132.0967866666666
$ws.Range("N15").Value = 396.29036
$ws.Range("O15").Value = 0.6000298026593724
$ws.Range("P15").Value = 0.6000298026593726
$ws.Range("Q15").Value = 11592.06597779603
$ws.Range("R15").Value = 104328.5938001643
$ws.Range("S15").Value = 0.2391784137249169
$ws.Range("T15").Value = 0.2391784137249169
$ws.Range("G16").Value = 87.75433733333334
$ws.Range("H16").Value = 263.263012
$ws.Range("I16").Value = 0.3986108901005617
$ws.Range("J16").Value = 0.3986108901005617
$ws.Range("O16").Value = 0.0003198253300583355
$ws.Range("P16").Value = 0.0003198253300583356
$ws.Range("Q16").Value = 6.178753640194222
$ws.Range("R16").Value = 55.608782761748
$ws.Range("S16").Value = 0.0001274858594912591
$ws.Range("T16").Value = 0.0001274858594912591
$ws.Range("G17").Value = 87.75433733333334
$ws.Range("H17").Value = 263.263012
$ws.Range("I17").Value = 0.3986108901005617
$ws.Range("J17").Value = 0.3986108901005617
$ws.Range("M17").Value = 87.75433733333334
$ws.Range("N17").Value = 263.263012
$ws.Range("O17").Value = 0.3986108901005617
$ws.Range("P17").Value = 0.3986108901005617
$ws.Range("Q17").Value = 7700.823720812461
$ws.Range("R17").Value = 69307.41348731215
$ws.Range("S17").Value = 0.1588906417067621
$ws.Range("T17").Value = 0.1588906417067621
